$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for rows 2-5 from
# 2023-09-16 (serial 45185) to 2023-10-05 (serial 45204).
foreach ($row in 2..5) {
    $ws.Cells.Item($row, 3).Value = 45204
}
